$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# 1. Shape restructuring on the footer banner:
#    - "Rectangle 29" is replaced by a new "Terminator 32" flow-chart-terminator
#      shape (same footer-bar look, new geometry).
#    - "Rectangle 5" (the big footer background bar) and the new "Terminator 32"
#      shape move to the very front of the z-order (right after the group
#      shape properties, before "Picture 3").
# ---------------------------------------------------------------------------

$rect29 = $s.Shapes.Item("Rectangle 29")

# Duplicate Rectangle 29 so the new shape inherits identical style/fill/line
# formatting (p:style block, schemeClr references, etc.), then re-purpose it.
$dupRange = $rect29.Duplicate()
$term32 = $dupRange.Item(1)
$term32.Name = "Terminator 32"
$term32.AutoShapeType = 69   # msoShapeFlowchartTerminator
$term32.Left = 21031200 / 12700.0
$term32.Top = 31242000 / 12700.0
$term32.Width = 14020800 / 12700.0
$term32.Height = 1676400 / 12700.0

# Remove the original Rectangle 29 shape entirely.
$rect29.Delete()

# Send "Rectangle 5" all the way to the back (front of the XML / bottom of
# the stack) so it becomes the first shape in the tree.
$rect5 = $s.Shapes.Item("Rectangle 5")
$rect5.ZOrder(1)    # msoSendToBack

# Put the new Terminator shape right after Rectangle 5 (second shape).
$term32.ZOrder(1)   # msoSendToBack
$term32.ZOrder(2)   # msoBringForward (now sits right after Rectangle 5)

# ---------------------------------------------------------------------------
# 2. Collapse the "Early Adopters" run fragments into single runs per
#    paragraph (text content unchanged, just merged run boundaries).
# ---------------------------------------------------------------------------

$early = $s.Shapes.Item("Rectangle 21")
$tr = $early.TextFrame.TextRange

function Merge-Segment($range, $marker, $newText) {
    $full = $range.Text
    $idx = $full.IndexOf($marker)
    if ($idx -ge 0) {
        $sub = $range.Characters($idx + 1, $newText.Length)
        $sub.Text = $newText
    }
}

Merge-Segment $tr "XML Data Service Style" "XML Data Service Style: transfer large xml document and large number of xml documents (via Zip stream)"
Merge-Segment $tr "IVI Middleware CERR Data Service" "IVI Middleware CERR Data Service: transfer CERR objects both in upload and download"
Merge-Segment $tr "NCIA Data Service" "NCIA Data Service: transfer DICOM images, both upload and download"
Merge-Segment $tr "IVI Middleware DICOM Data Service" "IVI Middleware DICOM Data Service: transfer DICOM images, both upload and download"
